# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the first-of-quarter date for each observation
# (e.g. 1988-07-01, 1988-10-01, ...). The fix re-indexes each date to the
# 15th of the quarter's middle month (e.g. 1988-07-01 -> 1988-08-15),
# leaving column B (the revision values) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    $d = [datetime]::FromOADate($serial)
    $newDate = $d.AddMonths(1).AddDays(14)
    $cell.Value2 = $newDate.ToOADate()
}
